$d = $word.ActiveDocument

# The document contains three "<id>...</id>" tags (p093v_1, p093v_2, p093v_3)
# that are each split across three separate runs:
#   run1: "<id>"     (Courier New, color 7f6000, sz 18)
#   run2: "p093v_N"  (plain, color 000000)
#   run3: "</id>"    (Courier New, color 7f6000, sz 18)
#
# Collapse each triple into a single run "<id>p093v_N</id>" that keeps the
# Courier New / 7f6000 / sz18 formatting of the surrounding tag runs.
# Using Find & Replace across the whole matched range merges the runs and
# Word applies the formatting of the first run in the match to the new text.

$ids = @("p093v_1", "p093v_2", "p093v_3")

foreach ($id in $ids) {
    $target = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, $target, 2) | Out-Null
}
